# Add a new "24-jul" column (AG) to the right of the existing data (A:AF)
# with the new day's values for each row, matching the formatting used by
# the preceding columns (text header, centered integer values).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column (text-formatted, same as the other date headers)
$ws.Range("AG1").Value = "24-jul"
$ws.Range("AG1").NumberFormat = "@"

# New values for the new column, one per data row (rows 2-11)
$values = @(12, 15, 10, 17, 18, 12, 18, 16, 27, 21)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 2 + $i
    $cell = $ws.Cells.Item($row, 33)
    $cell.Value = $values[$i]
    $cell.NumberFormat = "0"
    $cell.HorizontalAlignment = -4108
}

# Update the active selection to mirror the post-edit state (AG12)
$ws.Range("AG12").Select()
